$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the header date by one day (45308 -> 45309, i.e. 2024-01-17 -> 2024-01-18)
$ws.Range("A1").Value = 45309

# Step 2: update prices for CAJA p/ ENROLLADOR CHICA (D30) and GRANDE (D31)
$ws.Range("D30").Value = 570
$ws.Range("D31").Value = 690
